$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Sending cluster: MuSCs, Target cluster: ECs)
$ws.Range("G2").Value = 0.01135533333333333
$ws.Range("H2").Value = 0.034066
$ws.Range("M2").Value = 0.008446
$ws.Range("N2").Value = 0.025338
$ws.Range("O2").Value = 0.001125187475737063
$ws.Range("P2").Value = 0.001125187475737063
$ws.Range("Q2").Value = 0.00009590714533333334
$ws.Range("R2").Value = 0.000863164308
$ws.Range("S2").Value = 0.001125187475737063
$ws.Range("T2").Value = 0.001125187475737063

# Row 3 (Sending cluster: MuSCs, Target cluster: FAPs)
$ws.Range("G3").Value = 0.01135533333333333
$ws.Range("H3").Value = 0.034066
$ws.Range("N3").Value = 9.970262999999999
$ws.Range("O3").Value = 0.4427506139949732
$ws.Range("P3").Value = 0.4427506139949733
$ws.Range("Q3").Value = 0.037738553262
$ws.Range("R3").Value = 0.339646979358
$ws.Range("S3").Value = 0.4427506139949732
$ws.Range("T3").Value = 0.4427506139949733

# Row 4 (Sending cluster: MuSCs, Target cluster: MuSCs)
$ws.Range("G4").Value = 0.01135533333333333
$ws.Range("H4").Value = 0.034066
$ws.Range("M4").Value = 4.174437666666667
$ws.Range("N4").Value = 12.523313
$ws.Range("O4").Value = 0.5561241985292896
$ws.Range("P4").Value = 0.5561241985292896
$ws.Range("Q4").Value = 0.04740213118422223
$ws.Range("R4").Value = 0.426619180658
$ws.Range("S4").Value = 0.5561241985292896
$ws.Range("T4").Value = 0.5561241985292896
